# Auto-generated edit script: applies scheduled market-data refresh to Behemoth_Profits sheets
# (mirrors the upstream commit "chore: update Sheets via scheduled runner")

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")

# updated values
$ws.Range("H11").Value = 1308.4
$ws.Range("I11").Value = 1308.4
$ws.Range("K11").Value = 1308.4
$ws.Range("M11").Value = -1168.4
$ws.Range("H116").Value = 7013.357
$ws.Range("I116").Value = 6900
$ws.Range("J116").Value = 7098.375
$ws.Range("K116").Value = 6900
$ws.Range("L116").Value = 7098.375
$ws.Range("M116").Value = -3458
$ws.Range("N116").Value = -13982.375
$ws.Range("H137").Value = 5930.909
$ws.Range("I137").Value = 1885.2858
$ws.Range("J137").Value = 13010.75
$ws.Range("K137").Value = 5655.857400000001
$ws.Range("L137").Value = 39032.25
$ws.Range("M137").Value = -3105.857400000001
$ws.Range("N137").Value = -44132.25
$ws.Range("H138").Value = 3088.8096
$ws.Range("J138").Value = 3193.45
$ws.Range("L138").Value = 9580.349999999999
$ws.Range("N138").Value = -19860.35

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")

# updated values
$ws.Range("H97").Value = 1568
$ws.Range("I97").Value = 1568
$ws.Range("K97").Value = 1568
$ws.Range("M97").Value = -1072
$ws.Range("H122").Value = 1278.2
$ws.Range("I122").Value = 1186.8889
$ws.Range("K122").Value = 3560.6667
$ws.Range("M122").Value = -1110.6667

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")

# updated values
$ws.Range("H80").Value = 1046.5714
$ws.Range("J80").Value = 1154.6666
$ws.Range("L80").Value = 1154.6666
$ws.Range("N80").Value = -3150.6666
$ws.Range("H83").Value = 1046.5714
$ws.Range("J83").Value = 1154.6666
$ws.Range("L83").Value = 5773.333000000001
$ws.Range("N83").Value = -15757.333
$ws.Range("H107").Value = 2316.7144
$ws.Range("I107").Value = 2316.7144
$ws.Range("K107").Value = 2316.7144
$ws.Range("M107").Value = -396.7143999999998

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")

# updated values
$ws.Range("H16").Value = 2201.125
$ws.Range("I16").Value = 1934.8334
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 1934.8334
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -1647.8334
$ws.Range("N16").Value = -3574
$ws.Range("H31").Value = 759338.8
$ws.Range("I31").Value = 9848.111000000001
$ws.Range("J31").Value = 1797095.2
$ws.Range("K31").Value = 9848.111000000001
$ws.Range("L31").Value = 1797095.2
$ws.Range("M31").Value = -9553.111000000001
$ws.Range("N31").Value = -1797685.2
$ws.Range("H34").Value = 759338.8
$ws.Range("I34").Value = 9848.111000000001
$ws.Range("J34").Value = 1797095.2
$ws.Range("K34").Value = 9848.111000000001
$ws.Range("L34").Value = 1797095.2
$ws.Range("M34").Value = -9646.111000000001
$ws.Range("N34").Value = -1797499.2
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("H58").Value = 2539.7354
$ws.Range("I58").Value = 2500.9583
$ws.Range("J58").Value = 2632.8
$ws.Range("K58").Value = 2500.9583
$ws.Range("L58").Value = 2632.8
$ws.Range("M58").Value = -2297.9583
$ws.Range("N58").Value = -3038.8
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("H113").Value = 2201.125
$ws.Range("I113").Value = 1934.8334
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1934.8334
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 235.1666
$ws.Range("N113").Value = -7340
$ws.Range("H132").Value = 2435.2222
$ws.Range("I132").Value = 2435.2222
$ws.Range("K132").Value = 7305.6666
$ws.Range("M132").Value = -4775.6666
$ws.Range("H136").Value = 2539.7354
$ws.Range("I136").Value = 2500.9583
$ws.Range("J136").Value = 2632.8
$ws.Range("K136").Value = 7502.874899999999
$ws.Range("L136").Value = 7898.400000000001
$ws.Range("M136").Value = -4952.874899999999
$ws.Range("N136").Value = -12998.4

# cells cleared (no longer populated)
$ws.Range("N43").ClearContents()
$ws.Range("N101").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")

# updated values
$ws.Range("H36").Value = 3704.9
$ws.Range("I36").Value = 783.3333
$ws.Range("J36").Value = 29999
$ws.Range("K36").Value = 2349.9999
$ws.Range("L36").Value = 89997
$ws.Range("M36").Value = -2180.9999
$ws.Range("N36").Value = -90335
$ws.Range("H141").Value = 187434.23
$ws.Range("I141").Value = 384048
$ws.Range("J141").Value = 12666.444
$ws.Range("K141").Value = 1152144
$ws.Range("L141").Value = 37999.33199999999
$ws.Range("M141").Value = -1146964
$ws.Range("N141").Value = -48359.33199999999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")

# updated values
$ws.Range("H80").Value = 1955.4286
$ws.Range("I80").Value = 1955.4286
$ws.Range("K80").Value = 1955.4286
$ws.Range("M80").Value = -957.4286
$ws.Range("H83").Value = 1955.4286
$ws.Range("I83").Value = 1955.4286
$ws.Range("K83").Value = 9777.143
$ws.Range("M83").Value = -4785.143
$ws.Range("H93").Value = 59982.668
$ws.Range("J93").Value = 59982.668
$ws.Range("L93").Value = 59982.668
$ws.Range("N93").Value = -63726.668
$ws.Range("H113").Value = 4132.727
$ws.Range("I113").Value = 3688.75
$ws.Range("K113").Value = 3688.75
$ws.Range("M113").Value = -1518.75
$ws.Range("H122").Value = 1591.2084
$ws.Range("I122").Value = 1551.6957
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4655.0871
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2205.0871
$ws.Range("N122").Value = -12400

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")

# updated values
$ws.Range("H16").Value = 1636.4615
$ws.Range("I16").Value = 1455.0714
$ws.Range("J16").Value = 1848.0834
$ws.Range("K16").Value = 1455.0714
$ws.Range("L16").Value = 1848.0834
$ws.Range("M16").Value = -1285.0714
$ws.Range("N16").Value = -2188.0834
$ws.Range("H40").Value = 2797.4092
$ws.Range("I40").Value = 1971.625
$ws.Range("K40").Value = 1971.625
$ws.Range("M40").Value = -1835.625
$ws.Range("H46").Value = 3181.8965
$ws.Range("I46").Value = 2851.238
$ws.Range("J46").Value = 4049.875
$ws.Range("K46").Value = 2851.238
$ws.Range("L46").Value = 4049.875
$ws.Range("M46").Value = -2663.238
$ws.Range("N46").Value = -4425.875
$ws.Range("H50").Value = 50001
$ws.Range("J50").Value = 50001
$ws.Range("L50").Value = 50001

# newly populated cells
$ws.Range("N50").Value = -51275

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")

# updated values
$ws.Range("H81").Value = 1550
$ws.Range("I81").Value = 1400
$ws.Range("J81").Value = 1700
$ws.Range("K81").Value = 2800
$ws.Range("L81").Value = 3400
$ws.Range("H84").Value = 1550
$ws.Range("I84").Value = 1400
$ws.Range("J84").Value = 1700
$ws.Range("K84").Value = 14000
$ws.Range("L84").Value = 17000
$ws.Range("H122").Value = 5304.4443
$ws.Range("I122").Value = 3245.4614
$ws.Range("J122").Value = 10657.8
$ws.Range("K122").Value = 9736.3842
$ws.Range("L122").Value = 31973.4
$ws.Range("M122").Value = -7286.3842
$ws.Range("N122").Value = -36873.39999999999
$ws.Range("H126").Value = 2600
$ws.Range("I126").Value = 1700
$ws.Range("K126").Value = 5100
$ws.Range("M126").Value = -2630
$ws.Range("H132").Value = 1390.16
$ws.Range("I132").Value = 1390.16
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4170.48
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1640.48

# newly populated cells
$ws.Range("M81").Value = -1739
$ws.Range("N81").Value = -5522
$ws.Range("M84").Value = -8696
$ws.Range("N84").Value = -27608

# cells cleared (no longer populated)
$ws.Range("N132").ClearContents()

